# Change the table style ("Table Design") of the two-column totals table on
# slide 16 from the deck's custom "Table_0" style to the built-in
# "Medium Style 2 - Accent 1" table style, just like picking a new style
# from the Table Design > Table Styles gallery in PowerPoint would do.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(16)

# The table lives in the 3rd shape (a graphicFrame) on this slide.
$tableShape = $slide.Shapes.Item(3)

if ($tableShape.HasTable) {
    $table = $tableShape.Table
    $table.ApplyStyle("{273B40A6-295E-4DEB-8C16-A31A7780890C}")
}
